$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the date series: copy formatting (date number format) from A47 into A48,
# then set the new date and activity values for 2020-12-01.
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A48").Value = 44166
$ws.Range("B48").Value = 9

# Move the active selection to B49, as Excel does after entering the new row.
$ws.Range("B49").Select()
